# Generate Report for Handoff
# Adds a new handoff row (bf20a7b4-dc56-4daa-b9bf-42b5d4fcf2ce) to the
# Overview, zh-cn and de-de sheets of the localization-status report,
# mirroring the existing 75079753-... row that is already present.

$wb = $excel.ActiveWorkbook

$newId        = "bf20a7b4-dc56-4daa-b9bf-42b5d4fcf2ce"
$newFile      = "$newId.md"
$newHash      = "b22d9ccd1c7c696d5ff9abe3db34e18ff8f9ca1b"
$statusText   = "Ready for handoff"
$extText      = ".md"
$includeText  = "Include"
$epochText    = "0001-01-01 00:00:00"

$mdHyperlinkTarget = "https://github.com/OpenLocalizationTest/oltest/blob/40982135b1f85376ec3dd62ee6a03d223cd854c4/e2e/$newFile"

# ---------------------------------------------------------------------
# Overview sheet — one summary row per handed-off file
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
$wsOverview.Range("D3").Value = "2016-03-22 06:39:06"
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdHyperlinkTarget, "", "", $newFile)

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$newId.$newHash.zh-cn.xlf"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ec94c1b760e569d4c81319d2da74300473bd6e8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Range("B3").Value = $extText
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("D3").Value = $zhXlf
$wsZhCn.Range("E3").Value = "2016-03-22 06:39:02"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = $epochText
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J3").Value = $includeText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdHyperlinkTarget, "", "", $newFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhXlfTarget, "", "", $zhXlf)

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf = "$newId.$newHash.de-de.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/699eaa0b16cf04956610af6161471796de7002a1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Range("B3").Value = $extText
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("D3").Value = $deXlf
$wsDeDe.Range("E3").Value = "2016-03-22 06:39:06"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = $epochText
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("J3").Value = $includeText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdHyperlinkTarget, "", "", $newFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deXlfTarget, "", "", $deXlf)
